$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# 1 & 8: Title change (appears twice, same old -> same new text)
Replace-Text "Play Buffalo Boost for Free - Slot Game Review" "Play Buffalo Boost Free - Slot Game Review"

# 2: bullet "What we like" item 1
Replace-Text "Collect Feature with chance to trigger Free Spins mode" "Beautiful Grand Canyon scenery"

# 3: bullet "What we like" item 2
Replace-Text "Buy Feature allows players to activate the bonus mode with a small fee" "Calming soundtrack"

# 4: bullet "What we like" item 3
Replace-Text "Wild feature with buffalo symbol as most lucrative and aids in replacing other symbols" "High-quality gameplay"

# 5: bullet "What we like" item 4
Replace-Text "Beautiful Grand Canyon scenery and calming soundtrack" "Original bonus features"

# 6: bullet "What we don't like" item 1
Replace-Text "Poker card symbols have minimal significance" "Limited significance of poker card symbols"

# 7: bullet "What we don't like" item 2
Replace-Text "Free Spins mode may take a while to trigger" "Small fee required to activate bonus mode"

# 9: meta description
Replace-Text "Read our review of Buffalo Boost slot game and play for free. Learn about Collect and Buy Feature, pay lines, symbols, and gameplay experience." "Read our review of Buffalo Boost, a slot game set in the Grand Canyon. Play for free and enjoy original bonus features."
